$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update F column "想去人数" values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 633
$ws1.Range("F5").Value = 556
$ws1.Range("F6").Value = 308
$ws1.Range("F7").Value = 2767
$ws1.Range("F9").Value = 7756
$ws1.Range("F10").Value = 201
$ws1.Range("F11").Value = 469
$ws1.Range("F13").Value = 326

# Sheet "全部类型" (fourth sheet) - update F column "想去人数" values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 633
$ws4.Range("F5").Value = 556
$ws4.Range("F6").Value = 308
$ws4.Range("F9").Value = 2767
$ws4.Range("F11").Value = 7756
$ws4.Range("F12").Value = 201
$ws4.Range("F13").Value = 469
$ws4.Range("F17").Value = 326
